# Updates the cryptos price/volume table (and the Kaspa / EthereumClassic row
# swap) to match the refreshed "cryptos list" data, per the commit diff.
# Numeric-looking Price values in column D are written with a leading
# apostrophe so Excel keeps them as literal text (preserving trailing/leading
# zeros such as "0.0830" or "0.130") instead of re-casting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.399.83'
$ws.Range('E2').Value = '  +1.32%  '
$ws.Range('D3').Value = '1.943.06'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''242.96'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').Value = '''0.618'
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('D7').Value = '''58.42'
$ws.Range('E7').Value = '  -2.85%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').Value = '''0.362'
$ws.Range('E9').Value = '  -2.40%  '
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').Value = '''0.0830'
$ws.Range('E11').Value = '  +3.51%  '
$ws.Range('E12').Value = '  +1.76%  '
$ws.Range('D13').Value = '''21.55'
$ws.Range('E13').Value = '  -2.17%  '
$ws.Range('D14').Value = '''0.821'
$ws.Range('E14').Value = '  -3.49%  '
$ws.Range('D15').Value = '2.228.61'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').Value = '''13.59'
$ws.Range('E16').Value = '  -2.38%  '
$ws.Range('D17').Value = '''5.22'
$ws.Range('E17').Value = '  -3.03%  '
$ws.Range('D18').Value = '1.980.66'
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('D19').Value = '36.310.58'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').Value = '''69.51'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('D22').Value = '''228.54'
$ws.Range('E22').Value = '  -2.65%  '
$ws.Range('E23').Value = '  -2.75%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '''2.44'
$ws.Range('E25').Value = '  -3.18%  '
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').Value = '''9.18'
$ws.Range('E27').Value = '  -5.15%  '
$ws.Range('D28').Value = '''161.61'
$ws.Range('E28').Value = '  +1.10%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''19.47'
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '''0.130'
$ws.Range('E30').Value = '  +2.23%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('E32').Value = '  +1.75%  '
$ws.Range('E33').Value = '  -3.57%  '
$ws.Range('D34').Value = '''0.0626'
$ws.Range('E35').Value = '  -2.61%  '
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('E38').Value = '  -2.17%  '
$ws.Range('D39').Value = '''2.13'
$ws.Range('E39').Value = '  -5.52%  '
$ws.Range('D40').Value = '''3.02'
$ws.Range('E40').Value = '  -0.70%  '
$ws.Range('D41').Value = '''0.0979'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('E43').Value = '  -3.92%  '
$ws.Range('E44').Value = '  -1.00%  '
$ws.Range('D45').Value = '''15.99'
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('D46').Value = '1.354.44'
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('E47').Value = '  -4.44%  '
$ws.Range('D48').Value = '''87.75'
$ws.Range('E48').Value = '  -4.17%  '
$ws.Range('D49').Value = '''7.10'
$ws.Range('E49').Value = '  -4.47%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = '''45.31'
$ws.Range('E51').Value = '  +3.74%  '
